$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 179.1580256666667
$ws.Range("H2").Value = 537.4740770000001
$ws.Range("I2").Value = 0.3468013736386751
$ws.Range("J2").Value = 0.3468013736386751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 20969.47282209779
$ws.Range("R2").Value = 188725.2553988801
$ws.Range("S2").Value = 0.1125497263794132
$ws.Range("T2").Value = 0.1125497263794131
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 179.1580256666667
$ws.Range("H3").Value = 537.4740770000001
$ws.Range("I3").Value = 0.3468013736386751
$ws.Range("J3").Value = 0.3468013736386751
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 18198.8789357863
$ws.Range("R3").Value = 163789.9104220766
$ws.Range("S3").Value = 0.09767908149204027
$ws.Range("T3").Value = 0.09767908149204024
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 179.1580256666667
$ws.Range("H4").Value = 537.4740770000001
$ws.Range("I4").Value = 0.3468013736386751
$ws.Range("J4").Value = 0.3468013736386751
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 25445.2391687356
$ws.Range("R4").Value = 229007.1525186204
$ws.Range("S4").Value = 0.1365725657672217
$ws.Range("T4").Value = 0.1365725657672217
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 239.807332
$ws.Range("H5").Value = 719.421996
$ws.Range("I5").Value = 0.4642019905988459
$ws.Range("J5").Value = 0.4642019905988459
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 28068.14437813592
$ws.Range("R5").Value = 252613.2994032233
$ws.Range("S5").Value = 0.15065051928287
$ws.Range("T5").Value = 0.15065051928287
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 239.807332
$ws.Range("H6").Value = 719.421996
$ws.Range("I6").Value = 0.4642019905988459
$ws.Range("J6").Value = 0.4642019905988459
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 24359.63773736706
$ws.Range("R6").Value = 219236.7396363036
$ws.Range("S6").Value = 0.1307458029728386
$ws.Range("T6").Value = 0.1307458029728385
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 239.807332
$ws.Range("H7").Value = 719.421996
$ws.Range("I7").Value = 0.4642019905988459
$ws.Range("J7").Value = 0.4642019905988459
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 34059.06542255273
$ws.Range("R7").Value = 306531.5888029746
$ws.Range("S7").Value = 0.1828056683431374
$ws.Range("T7").Value = 0.1828056683431374
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 97.63589966666666
$ws.Range("H8").Value = 292.907699
$ws.Range("I8").Value = 0.1889966357624789
$ws.Range("J8").Value = 0.1889966357624789
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 11427.75120959685
$ws.Range("R8").Value = 102849.7608863716
$ws.Range("S8").Value = 0.06133631888049831
$ws.Range("T8").Value = 0.06133631888049831
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 97.63589966666666
$ws.Range("H9").Value = 292.907699
$ws.Range("I9").Value = 0.1889966357624789
$ws.Range("J9").Value = 0.1889966357624789
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 9917.858333213588
$ws.Range("R9").Value = 89260.72499892229
$ws.Range("S9").Value = 0.05323225105099719
$ws.Range("T9").Value = 0.05323225105099718
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 97.63589966666666
$ws.Range("H10").Value = 292.907699
$ws.Range("I10").Value = 0.1889966357624789
$ws.Range("J10").Value = 0.1889966357624789
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 13866.91335332814
$ws.Range("R10").Value = 124802.2201799532
$ws.Range("S10").Value = 0.07442806583098345
$ws.Range("T10").Value = 0.07442806583098345